# DeveloperGuide: Revise some of the figures to match color scheme
#
# 1) Refresh the cached "last saved" date/time field text (datetimeFigureOut)
#    on the slide master, every slide layout, and the notes master.
# 2) On the UndoRedo sequence diagram slide: remove the stray "Logic" label
#    rectangle, and rename the ":Address" lifeline to ":Coin".

function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame) {
                if ($shp.TextFrame.TextRange.Text -eq "7/29/2017") {
                    $shp.TextFrame.TextRange.Text = "08-Apr-18"
                }
            }
        }
    }
}

$p = $ppt.ActivePresentation

# --- Slide master date placeholder ---
$m = $p.SlideMaster
Update-DatePlaceholder $m.Shapes

# --- Every slide layout's date placeholder ---
$cls = $m.CustomLayouts
for ($i = 1; $i -le $cls.Count; $i++) {
    $cl = $cls.Item($i)
    Update-DatePlaceholder $cl.Shapes
}

# --- Notes master date placeholder ---
# (direct TextFrame.TextRange write is a no-op for the notes master in this
# host, so go through HeadersFooters.DateAndTime instead)
$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "08-Apr-18"

# --- Slide content edits ---
$s = $p.Slides.Item(1)

# Remove the "Logic" rounded-rectangle label (shape 1, "Rectangle 65").
$s.Shapes.Item(1).Delete()

# After the delete, shape 4 is "Rectangle 62" (id 16) whose first paragraph
# reads ":Address" (second paragraph continues "BookParser"). Rename the
# lifeline from ":Address" to ":Coin".
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange
$chars = $tr.Characters(1, 8)
$chars.Text = ":Coin"

Write-Host "Edit complete"
